$wb = $excel.ActiveWorkbook
$new = $wb.Worksheets.Add()
$new.Name = "Temp"
$old = $wb.Worksheets.Item("Sheet1")
$old.Delete()
$new.Name = "Sheet1"
$new.Range("A1").Value = "For internal use only. Do not fill in or change"
$new.Range("A2").Value = "BusinessKey"
$new.Range("B2").Value = "Code"
$new.Range("C2").Value = "Framework_ID"
$new.Range("D2").Value = "Name"
$new.Range("E2").Value = "OrganizationBusinessKey"
$new.Range("A2:E2").Font.Bold = $true
$new.Range("A2:E2").Font.Underline = $true
